$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 25

$ws.Cells.Item($row, 1).Value = 1
$ws.Cells.Item($row, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item($row, 3).Value = "Arica y Parinacota"
$ws.Cells.Item($row, 4).Value = Get-Date -Year 2022 -Month 5 -Day 25 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($row, 5).Value = 15
$ws.Cells.Item($row, 6).Value = 100112013
$ws.Cells.Item($row, 7).Value = "Alcachofa"
$ws.Cells.Item($row, 8).Value = "Madrigal"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 250
$ws.Cells.Item($row, 11).Value = 21000
$ws.Cells.Item($row, 12).Value = 22000
$ws.Cells.Item($row, 13).Value = 21500
$ws.Cells.Item($row, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item($row, 15).Value = "Región de Coquimbo"
$ws.Cells.Item($row, 16).Value = 538
$ws.Cells.Item($row, 17).Value = 40
$ws.Cells.Item($row, 18).Value = "Hortaliza"
